# Applies:
#  1) slide16's table: tableStyleId {935A1E24-...} -> {03939F60-EB60-47FA-A3B3-BAA6A5005200}
#  2) theme color swap between the "Office Theme" (originally theme1.xml, only used
#     by the Notes Master) and the "Integral" theme (originally theme2.xml, used by
#     the Slide Master / live presentation design) - i.e. the live design reverts to
#     the default Office Theme colors.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{03939F60-EB60-47FA-A3B3-BAA6A5005200}")
    }
}

# --- 2) Theme colours: restore the default Office Theme palette -----------------
function ToOleRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    @(0x00, 0x00, 0x00), # 1  dk1
    @(0xFF, 0xFF, 0xFF), # 2  lt1
    @(0x44, 0x54, 0x6A), # 3  dk2
    @(0xE7, 0xE6, 0xE6), # 4  lt2
    @(0x5B, 0x9B, 0xD5), # 5  accent1
    @(0xED, 0x7D, 0x31), # 6  accent2
    @(0xA5, 0xA5, 0xA5), # 7  accent3
    @(0xFF, 0xC0, 0x00), # 8  accent4
    @(0x44, 0x72, 0xC4), # 9  accent5
    @(0x70, 0xAD, 0x47), # 10 accent6
    @(0x05, 0x63, 0xC1), # 11 hlink
    @(0x95, 0x4F, 0x72)  # 12 folHlink
)

$cs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $c = $officeColors[$i - 1]
    $cs.Colors($i).RGB = ToOleRgb $c[0] $c[1] $c[2]
}
